# "Legs and Room 3 Sesi 2 Update!"
# Updates the calibration measurements for LEG rows 1-5 (sheet rows 4-8).
# Dependent formulas in columns H, I, J, K, L, M recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- LEG L1 (row 4) ---
$ws.Range("C4").Value = 1500
$ws.Range("D4").Value = 1550
$ws.Range("F4").Value = 2000
$ws.Range("G4").Value = 2100

# --- LEG L2 (row 5) ---
$ws.Range("D5").Value = 1450

# --- LEG L3 (row 6) ---
$ws.Range("C6").Value = 1350
$ws.Range("D6").Value = 1650

# --- LEG R1 (row 7) ---
$ws.Range("D7").Value = 1350
$ws.Range("G7").Value = 850

# --- LEG R2 / Room 3 (row 8) ---
$ws.Range("C8").Value = 1550
$ws.Range("D8").Value = 1300
$ws.Range("E8").Value = 1870

# Match the author's final on-screen selection/scroll state.
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F14").Select()
